# Update KWD-loginGmail.xlsx to the latest "olo" placeholder convention:
# replace the <<username>>/<<password>> markers with {{username}}/{{password}}
# on Sheet1, and make Sheet1 (cell C3) the active selection instead of
# dataProvider.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("C2").Value = "{{username}}"
$ws1.Range("C3").Value = "{{password}}"

# Activate Sheet1 and select C3 so it becomes the saved selection/active tab.
$ws1.Activate()
[void]$ws1.Range("C3").Select()
